$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("İş Takip Listesi")

# Row 2
$ws.Range("G2").Value = "'3"
$ws.Range("G2").Style = "Normal"
$ws.Range("I2").Value = "İhaleli"
$ws.Range("J2").Value = "'2025-02-04"
$ws.Range("J2").Style = "Normal"
$ws.Range("L2").Value = "YER TESLİMİ YAPILDI"

# Row 3
$ws.Range("G3").Value = "'5"
$ws.Range("G3").Style = "Normal"
$ws.Range("I3").Value = "İhaleli"
$ws.Range("J3").Value = "'2025-07-08"
$ws.Range("J3").Style = "Normal"
$ws.Range("L3").Value = "YER TESLİMİ YAPILDI"
